$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Schedule")
$ws1.Range("A2").Value = 46040
$ws1.Range("B2").Value = 46040.16666666666
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 15.12
$ws1.Range("E2").Value = 363.42217275
$ws1.Range("F2").Value = 24.03585798611111

$ws1.Range("A3").Value = 46040.29166666666
$ws1.Range("B3").Value = 46040.79166666666
$ws1.Range("C3").Value = 12
$ws1.Range("D3").Value = 45.36
$ws1.Range("E3").Value = -62.45913375000001
$ws1.Range("F3").Value = -1.376965029761905

$ws2 = $wb.Worksheets.Item("Detailed")
$ws2.Range("A2").Value = 46040
$ws2.Range("B2").Value = 56.98
$ws2.Range("C2").Value = "historical"
$ws2.Range("D2").Value = 46040
$ws2.Range("E2").Value = "ON"

$ws2.Range("A3").Value = 46040.02083333334
$ws2.Range("B3").Value = 36.2
$ws2.Range("C3").Value = "historical"
$ws2.Range("D3").Value = 46040
$ws2.Range("E3").Value = "ON"

$ws2.Range("A4").Value = 46040.04166666666
$ws2.Range("B4").Value = 36.0603
$ws2.Range("C4").Value = "historical"
$ws2.Range("D4").Value = 46040
$ws2.Range("E4").Value = "ON"

$ws2.Range("A5").Value = 46040.0625
$ws2.Range("B5").Value = 36.2
$ws2.Range("C5").Value = "historical"
$ws2.Range("D5").Value = 46040
$ws2.Range("E5").Value = "ON"

$ws2.Range("A6").Value = 46040.08333333334
$ws2.Range("B6").Value = 36.2
$ws2.Range("C6").Value = "historical"
$ws2.Range("D6").Value = 46040
$ws2.Range("E6").Value = "ON"

$ws2.Range("A7").Value = 46040.10416666666
$ws2.Range("B7").Value = 56.97996
$ws2.Range("C7").Value = "forecast"
$ws2.Range("D7").Value = 46040
$ws2.Range("E7").Value = "ON"

$ws2.Range("A8").Value = 46040.125
$ws2.Range("B8").Value = 57.06021
$ws2.Range("C8").Value = "forecast"
$ws2.Range("D8").Value = 46040
$ws2.Range("E8").Value = "ON"

$ws2.Range("A9").Value = 46040.14583333334
$ws2.Range("B9").Value = 57.06022
$ws2.Range("C9").Value = "forecast"
$ws2.Range("D9").Value = 46040
$ws2.Range("E9").Value = "ON"

$ws2.Range("A10").Value = 46040.16666666666
$ws2.Range("B10").Value = 57.06022
$ws2.Range("C10").Value = "forecast"
$ws2.Range("D10").Value = 46040
$ws2.Range("E10").Value = "OFF"

$ws2.Range("A11").Value = 46040.1875
$ws2.Range("B11").Value = 57.06003
$ws2.Range("C11").Value = "forecast"
$ws2.Range("D11").Value = 46040
$ws2.Range("E11").Value = "OFF"

$ws2.Range("A12").Value = 46040.20833333334
$ws2.Range("B12").Value = 57.06003
$ws2.Range("C12").Value = "forecast"
$ws2.Range("D12").Value = 46040
$ws2.Range("E12").Value = "OFF"

$ws2.Range("A13").Value = 46040.22916666666
$ws2.Range("B13").Value = 57.06003
$ws2.Range("C13").Value = "forecast"
$ws2.Range("D13").Value = 46040
$ws2.Range("E13").Value = "OFF"

$ws2.Range("A14").Value = 46040.25
$ws2.Range("B14").Value = 57.06003
$ws2.Range("C14").Value = "forecast"
$ws2.Range("D14").Value = 46040
$ws2.Range("E14").Value = "OFF"

$ws2.Range("A15").Value = 46040.27083333334
$ws2.Range("B15").Value = 51.00032
$ws2.Range("C15").Value = "forecast"
$ws2.Range("D15").Value = 46040
$ws2.Range("E15").Value = "OFF"

$ws2.Range("A16").Value = 46040.29166666666
$ws2.Range("B16").Value = 24.52381
$ws2.Range("C16").Value = "forecast"
$ws2.Range("D16").Value = 46040
$ws2.Range("E16").Value = "ON"

$ws2.Range("A17").Value = 46040.3125
$ws2.Range("B17").Value = 0.7
$ws2.Range("C17").Value = "forecast"
$ws2.Range("D17").Value = 46040
$ws2.Range("E17").Value = "ON"

$ws2.Range("A18").Value = 46040.33333333334
$ws2.Range("B18").Value = -5.51
$ws2.Range("C18").Value = "forecast"
$ws2.Range("D18").Value = 46040
$ws2.Range("E18").Value = "ON"

$ws2.Range("A19").Value = 46040.35416666666
$ws2.Range("B19").Value = 0.6145699999999999
$ws2.Range("C19").Value = "forecast"
$ws2.Range("D19").Value = 46040
$ws2.Range("E19").Value = "ON"

$ws2.Range("A20").Value = 46040.375
$ws2.Range("B20").Value = 0.00961
$ws2.Range("C20").Value = "forecast"
$ws2.Range("D20").Value = 46040
$ws2.Range("E20").Value = "ON"

$ws2.Range("A21").Value = 46040.39583333334
$ws2.Range("B21").Value = -1.74343
$ws2.Range("C21").Value = "forecast"
$ws2.Range("D21").Value = 46040
$ws2.Range("E21").Value = "ON"

$ws2.Range("A22").Value = 46040.41666666666
$ws2.Range("B22").Value = 0
$ws2.Range("C22").Value = "forecast"
$ws2.Range("D22").Value = 46040
$ws2.Range("E22").Value = "ON"

$ws2.Range("A23").Value = 46040.4375
$ws2.Range("B23").Value = 0.7
$ws2.Range("C23").Value = "forecast"
$ws2.Range("D23").Value = 46040
$ws2.Range("E23").Value = "ON"

$ws2.Range("A24").Value = 46040.45833333334
$ws2.Range("B24").Value = 0
$ws2.Range("C24").Value = "forecast"
$ws2.Range("D24").Value = 46040
$ws2.Range("E24").Value = "ON"

$ws2.Range("A25").Value = 46040.47916666666
$ws2.Range("B25").Value = -0.91939
$ws2.Range("C25").Value = "forecast"
$ws2.Range("D25").Value = 46040
$ws2.Range("E25").Value = "ON"

$ws2.Range("A26").Value = 46040.5
$ws2.Range("B26").Value = -0.90235
$ws2.Range("C26").Value = "forecast"
$ws2.Range("D26").Value = 46040
$ws2.Range("E26").Value = "ON"

$ws2.Range("A27").Value = 46040.52083333334
$ws2.Range("B27").Value = 0
$ws2.Range("C27").Value = "forecast"
$ws2.Range("D27").Value = 46040
$ws2.Range("E27").Value = "ON"

$ws2.Range("A28").Value = 46040.54166666666
$ws2.Range("B28").Value = -5.51011
$ws2.Range("C28").Value = "forecast"
$ws2.Range("D28").Value = 46040
$ws2.Range("E28").Value = "ON"

$ws2.Range("A29").Value = 46040.5625
$ws2.Range("B29").Value = -6.8
$ws2.Range("C29").Value = "forecast"
$ws2.Range("D29").Value = 46040
$ws2.Range("E29").Value = "ON"

$ws2.Range("A30").Value = 46040.58333333334
$ws2.Range("B30").Value = -14
$ws2.Range("C30").Value = "forecast"
$ws2.Range("D30").Value = 46040
$ws2.Range("E30").Value = "ON"

$ws2.Range("A31").Value = 46040.60416666666
$ws2.Range("B31").Value = -23.5
$ws2.Range("C31").Value = "forecast"
$ws2.Range("D31").Value = 46040
$ws2.Range("E31").Value = "ON"

$ws2.Range("A32").Value = 46040.625
$ws2.Range("B32").Value = -14.7657
$ws2.Range("C32").Value = "forecast"
$ws2.Range("D32").Value = 46040
$ws2.Range("E32").Value = "ON"

$ws2.Range("A33").Value = 46040.64583333334
$ws2.Range("B33").Value = -14
$ws2.Range("C33").Value = "forecast"
$ws2.Range("D33").Value = 46040
$ws2.Range("E33").Value = "ON"

$ws2.Range("A34").Value = 46040.66666666666
$ws2.Range("B34").Value = -7.15426
$ws2.Range("C34").Value = "forecast"
$ws2.Range("D34").Value = 46040
$ws2.Range("E34").Value = "ON"

$ws2.Range("A35").Value = 46040.6875
$ws2.Range("B35").Value = -6.48482
$ws2.Range("C35").Value = "forecast"
$ws2.Range("D35").Value = 46040
$ws2.Range("E35").Value = "ON"

$ws2.Range("A36").Value = 46040.70833333334
$ws2.Range("B36").Value = -6
$ws2.Range("C36").Value = "forecast"
$ws2.Range("D36").Value = 46040
$ws2.Range("E36").Value = "ON"

$ws2.Range("A37").Value = 46040.72916666666
$ws2.Range("B37").Value = 0.66143
$ws2.Range("C37").Value = "forecast"
$ws2.Range("D37").Value = 46040
$ws2.Range("E37").Value = "ON"

$ws2.Range("A38").Value = 46040.75
$ws2.Range("B38").Value = 3.98333
$ws2.Range("C38").Value = "forecast"
$ws2.Range("D38").Value = 46040
$ws2.Range("E38").Value = "ON"

$ws2.Range("A39").Value = 46040.77083333334
$ws2.Range("B39").Value = 12.03666
$ws2.Range("C39").Value = "forecast"
$ws2.Range("D39").Value = 46040
$ws2.Range("E39").Value = "ON"

$ws2.Range("A40").Value = 46040.79166666666
$ws2.Range("B40").Value = 36.2
$ws2.Range("C40").Value = "forecast"
$ws2.Range("D40").Value = 46040
$ws2.Range("E40").Value = "OFF"

$ws2.Range("A41").Value = 46040.8125
$ws2.Range("B41").Value = 56.98
$ws2.Range("C41").Value = "forecast"
$ws2.Range("D41").Value = 46040
$ws2.Range("E41").Value = "OFF"

$ws2.Range("A42").Value = 46040.83333333334
$ws2.Range("B42").Value = 57.3
$ws2.Range("C42").Value = "forecast"
$ws2.Range("D42").Value = 46040
$ws2.Range("E42").Value = "OFF"

$ws2.Range("A43").Value = 46040.85416666666
$ws2.Range("B43").Value = 56.98
$ws2.Range("C43").Value = "forecast"
$ws2.Range("D43").Value = 46040
$ws2.Range("E43").Value = "OFF"

$ws2.Range("A44").Value = 46040.875
$ws2.Range("B44").Value = 56.98
$ws2.Range("C44").Value = "forecast"
$ws2.Range("D44").Value = 46040
$ws2.Range("E44").Value = "OFF"

$ws2.Range("A45").Value = 46040.89583333334
$ws2.Range("B45").Value = 46.20513
$ws2.Range("C45").Value = "forecast"
$ws2.Range("D45").Value = 46040
$ws2.Range("E45").Value = "OFF"

$ws2.Range("A46").Value = 46040.91666666666
$ws2.Range("B46").Value = 56.34431
$ws2.Range("C46").Value = "forecast"
$ws2.Range("D46").Value = 46040
$ws2.Range("E46").Value = "OFF"

$ws2.Range("A47").Value = 46040.9375
$ws2.Range("B47").Value = 36.0601
$ws2.Range("C47").Value = "forecast"
$ws2.Range("D47").Value = 46040
$ws2.Range("E47").Value = "OFF"

$ws2.Range("A48").Value = 46040.95833333334
$ws2.Range("B48").Value = 36.0601
$ws2.Range("C48").Value = "forecast"
$ws2.Range("D48").Value = 46040
$ws2.Range("E48").Value = "OFF"

$ws2.Range("A49").Value = 46040.97916666666
$ws2.Range("B49").Value = 36.06
$ws2.Range("C49").Value = "forecast"
$ws2.Range("D49").Value = 46040
$ws2.Range("E49").Value = "OFF"
$ws2.Range("A49").NumberFormat = $ws2.Range("A48").NumberFormat
$ws2.Range("D49").NumberFormat = $ws2.Range("D48").NumberFormat
